$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates species-observation data among rows 3, 4 and 6
# (row 3 -> gets what row 4 had, row 4 -> gets what row 6 had,
#  row 6 -> gets what row 3 had but with a new Taxonsorteringsordning
#  value in column B), and bumps column B on row 5.

# --- Row 3: becomes the old row 4 (Talltita) ---
$ws.Range("A3").Value = 112313576
$ws.Range("B3").Value = 56575
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 103021
$ws.Range("F3").Value = "Talltita"
$ws.Range("G3").Value = "Poecile montanus"
$ws.Range("H3").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"
$ws.Range("M3").Value = "lockläte, övriga läten"
$ws.Range("Q3").Value = 752714
$ws.Range("R3").Value = 7093570

# --- Row 4: becomes the old row 6 (Spillkråka) ---
$ws.Range("A4").Value = 112313590
$ws.Range("B4").Value = 56446
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = "Spillkråka"
$ws.Range("G4").Value = "Dryocopus martius"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("I4").Value = $null
$ws.Range("M4").Value = "äldre spår"
$ws.Range("Q4").Value = 752543
$ws.Range("R4").Value = 7093684

# --- Row 5: only the Taxonsorteringsordning (column B) changes ---
$ws.Range("B5").Value = 90837

# --- Row 6: becomes the old row 3 (Dropptaggsvamp), with a new column B ---
$ws.Range("A6").Value = 112313702
$ws.Range("B6").Value = 90814
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 4364
$ws.Range("F6").Value = "Dropptaggsvamp"
$ws.Range("G6").Value = "Hydnellum ferrugineum"
$ws.Range("H6").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I6").Value = $null
$ws.Range("M6").Value = $null
$ws.Range("Q6").Value = 752827
$ws.Range("R6").Value = 7093488
